$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 21.64429833926596
$ws.Cells.Item(2, 3).Value = 24.55962196639037
$ws.Cells.Item(2, 4).Value = 13.90695676658817
$ws.Cells.Item(2, 5).Value = 14.07553462178563
$ws.Cells.Item(2, 7).Value = 3.806446661780591
$ws.Cells.Item(2, 9).Value = 52.24426545471579
$ws.Cells.Item(2, 10).Value = 8.42747175141918
$ws.Cells.Item(2, 13).Value = 24.92264367281062
$ws.Cells.Item(3, 2).Value = 21.59696295488702
$ws.Cells.Item(3, 3).Value = 24.15985547819727
$ws.Cells.Item(3, 4).Value = 13.89956147134617
$ws.Cells.Item(3, 5).Value = 14.08350364267313
$ws.Cells.Item(3, 7).Value = 3.81260108385568
$ws.Cells.Item(3, 9).Value = 50.95250479852896
$ws.Cells.Item(3, 10).Value = 8.439053292572849
$ws.Cells.Item(3, 13).Value = 24.91609946912542
$ws.Cells.Item(4, 2).Value = 21.57983199231682
$ws.Cells.Item(4, 3).Value = 23.91942476517431
$ws.Cells.Item(4, 4).Value = 13.89778528670908
$ws.Cells.Item(4, 5).Value = 14.09040970326367
$ws.Cells.Item(4, 7).Value = 3.816560401324335
$ws.Cells.Item(4, 9).Value = 50.14348598000983
$ws.Cells.Item(4, 10).Value = 8.446622381130977
$ws.Cells.Item(4, 13).Value = 24.92007041850963
$ws.Cells.Item(5, 2).Value = 21.57584279358089
$ws.Cells.Item(5, 3).Value = 23.82283141623229
$ws.Cells.Item(5, 4).Value = 13.89775539602232
$ws.Cells.Item(5, 5).Value = 14.09372918813814
$ws.Cells.Item(5, 7).Value = 3.81821949181197
$ws.Cells.Item(5, 9).Value = 49.81012397382042
$ws.Cells.Item(5, 10).Value = 8.449822357228195
$ws.Cells.Item(5, 13).Value = 24.923690624004
$ws.Cells.Item(6, 2).Value = 21.5753607330834
$ws.Cells.Item(6, 3).Value = 23.80687921005852
$ws.Cells.Item(6, 4).Value = 13.89779230060132
$ws.Cells.Item(6, 5).Value = 14.09431086885514
$ws.Cells.Item(6, 7).Value = 3.818497746274611
$ws.Cells.Item(6, 9).Value = 49.75455660041102
$ws.Cells.Item(6, 10).Value = 8.450360697676242
$ws.Cells.Item(6, 13).Value = 24.92441244266952
$ws.Cells.Item(7, 2).Value = 21.57976609413287
$ws.Cells.Item(7, 3).Value = 23.91811631283733
$ws.Cells.Item(7, 4).Value = 13.8977820758213
$ws.Cells.Item(7, 5).Value = 14.09045242678051
$ws.Cells.Item(7, 7).Value = 3.816582591281122
$ws.Cells.Item(7, 9).Value = 50.13900462162998
$ws.Cells.Item(7, 10).Value = 8.446665068959893
$ws.Cells.Item(7, 13).Value = 24.92011114594413
$ws.Cells.Item(8, 2).Value = 21.62549652587349
$ws.Cells.Item(8, 3).Value = 24.42081171246445
$ws.Cells.Item(8, 4).Value = 13.90383242493083
$ws.Cells.Item(8, 5).Value = 14.0778639585577
$ws.Cells.Item(8, 7).Value = 3.808531402069028
$ws.Cells.Item(8, 9).Value = 51.8023420206065
$ws.Cells.Item(8, 10).Value = 8.431370224278897
$ws.Cells.Item(8, 13).Value = 24.91872535492563
$ws.Cells.Item(9, 2).Value = 21.81001551364464
$ws.Cells.Item(9, 3).Value = 25.44112318792602
$ws.Cells.Item(9, 4).Value = 13.93767716534248
$ws.Cells.Item(9, 5).Value = 14.06919354001799
$ws.Cells.Item(9, 7).Value = 3.794163409558882
$ws.Cells.Item(9, 9).Value = 54.92544661673688
$ws.Cells.Item(9, 10).Value = 8.40499519403858
$ws.Cells.Item(9, 13).Value = 24.97964447646673
$ws.Cells.Item(10, 2).Value = 22.00327332799591
$ws.Cells.Item(10, 3).Value = 26.20457972241439
$ws.Cells.Item(10, 4).Value = 13.97598705422258
$ws.Cells.Item(10, 5).Value = 14.07264374464818
$ws.Cells.Item(10, 7).Value = 3.784456748031787
$ws.Cells.Item(10, 9).Value = 57.11971336189745
$ws.Cells.Item(10, 10).Value = 8.387801479762796
$ws.Cells.Item(10, 13).Value = 25.06344233328516
$ws.Cells.Item(11, 2).Value = 22.10358714255183
$ws.Cells.Item(11, 3).Value = 26.55339096020074
$ws.Cells.Item(11, 4).Value = 13.99633711251616
$ws.Cells.Item(11, 5).Value = 14.07635624526034
$ws.Cells.Item(11, 7).Value = 3.780221840701191
$ws.Cells.Item(11, 9).Value = 58.09323033032612
$ws.Cells.Item(11, 10).Value = 8.380449245539259
$ws.Cells.Item(11, 13).Value = 25.11006228006177
$ws.Cells.Item(12, 2).Value = 22.14333798373451
$ws.Cells.Item(12, 3).Value = 26.68557161243073
$ws.Cells.Item(12, 4).Value = 14.00446293890428
$ws.Cells.Item(12, 5).Value = 14.07807086352363
$ws.Cells.Item(12, 7).Value = 3.778643894523932
$ws.Cells.Item(12, 9).Value = 58.45810929428402
$ws.Cells.Item(12, 10).Value = 8.37773226760091
$ws.Cells.Item(12, 13).Value = 25.12893728643234
$ws.Cells.Item(13, 2).Value = 22.13469881599357
$ws.Cells.Item(13, 3).Value = 26.65710191863227
$ws.Cells.Item(13, 4).Value = 14.00269424346469
$ws.Cells.Item(13, 5).Value = 14.07768784818975
$ws.Cells.Item(13, 7).Value = 3.778982593717545
$ws.Cells.Item(13, 9).Value = 58.37969714406717
$ws.Cells.Item(13, 10).Value = 8.378314435933991
$ws.Cells.Item(13, 13).Value = 25.12481794599091
$ws.Cells.Item(14, 2).Value = 22.10682223601598
$ws.Cells.Item(14, 3).Value = 26.56426436506996
$ws.Cells.Item(14, 4).Value = 13.99699722227541
$ws.Cells.Item(14, 5).Value = 14.07649111644317
$ws.Cells.Item(14, 7).Value = 3.780091508142882
$ws.Cells.Item(14, 9).Value = 58.12332592709785
$ws.Cells.Item(14, 10).Value = 8.380224374141173
$ws.Cells.Item(14, 13).Value = 25.11159067927727
$ws.Cells.Item(15, 2).Value = 22.08997613436142
$ws.Cells.Item(15, 3).Value = 26.5074072396037
$ws.Cells.Item(15, 4).Value = 13.99356226656519
$ws.Cells.Item(15, 5).Value = 14.07579831192677
$ws.Cells.Item(15, 7).Value = 3.780774092067044
$ws.Cells.Item(15, 9).Value = 57.96579377718978
$ws.Cells.Item(15, 10).Value = 8.381403002131812
$ws.Cells.Item(15, 13).Value = 25.10364754920628
$ws.Cells.Item(16, 2).Value = 21.99696559021021
$ws.Cells.Item(16, 3).Value = 26.18180381458997
$ws.Cells.Item(16, 4).Value = 13.97471590072793
$ws.Cells.Item(16, 5).Value = 14.07244430124237
$ws.Cells.Item(16, 7).Value = 3.784737122901197
$ws.Cells.Item(16, 9).Value = 57.05557581350135
$ws.Cells.Item(16, 10).Value = 8.388291380521386
$ws.Cells.Item(16, 13).Value = 25.06056658038256
$ws.Cells.Item(17, 2).Value = 21.94307055669186
$ws.Cells.Item(17, 3).Value = 25.98235431839935
$ws.Cells.Item(17, 4).Value = 13.96390253859009
$ws.Cells.Item(17, 5).Value = 14.07093606358748
$ws.Cells.Item(17, 7).Value = 3.787214409266845
$ws.Cells.Item(17, 9).Value = 56.49070352761542
$ws.Cells.Item(17, 10).Value = 8.392637136252658
$ws.Cells.Item(17, 13).Value = 25.03631445101849
$ws.Cells.Item(18, 2).Value = 21.91324012757073
$ws.Cells.Item(18, 3).Value = 25.86778566182647
$ws.Cells.Item(18, 4).Value = 13.95795813419529
$ws.Cells.Item(18, 5).Value = 14.07027024251174
$ws.Cells.Item(18, 7).Value = 3.788656303240897
$ws.Cells.Item(18, 9).Value = 56.16349437841159
$ws.Cells.Item(18, 10).Value = 8.395180889660516
$ws.Cells.Item(18, 13).Value = 25.02316549134165
$ws.Cells.Item(19, 2).Value = 21.90334125961231
$ws.Cells.Item(19, 3).Value = 25.82902418405977
$ws.Cells.Item(19, 4).Value = 13.95599273109718
$ws.Cells.Item(19, 5).Value = 14.07007942483484
$ws.Cells.Item(19, 7).Value = 3.789147435702562
$ws.Cells.Item(19, 9).Value = 56.05231745725657
$ws.Cells.Item(19, 10).Value = 8.396049760458149
$ws.Cells.Item(19, 13).Value = 25.01885093017886
$ws.Cells.Item(20, 2).Value = 21.94868695498543
$ws.Cells.Item(20, 3).Value = 26.00357150344896
$ws.Cells.Item(20, 4).Value = 13.96502516498216
$ws.Cells.Item(20, 5).Value = 14.07107574045426
$ws.Cells.Item(20, 7).Value = 3.786948937846628
$ws.Cells.Item(20, 9).Value = 56.55107586012556
$ws.Cells.Item(20, 10).Value = 8.392169952051651
$ws.Cells.Item(20, 13).Value = 25.03881331540762
$ws.Cells.Item(21, 2).Value = 22.11496255733475
$ws.Cells.Item(21, 3).Value = 26.59153143702326
$ws.Cells.Item(21, 4).Value = 13.99865919129496
$ws.Cells.Item(21, 5).Value = 14.07683424107847
$ws.Cells.Item(21, 7).Value = 3.779765097135802
$ws.Cells.Item(21, 9).Value = 58.19873239309737
$ws.Cells.Item(21, 10).Value = 8.379661559063843
$ws.Cells.Item(21, 13).Value = 25.1154427269055
$ws.Cells.Item(22, 2).Value = 22.23390431382199
$ws.Cells.Item(22, 3).Value = 26.97628726354567
$ws.Cells.Item(22, 4).Value = 14.02308679500212
$ws.Cells.Item(22, 5).Value = 14.08239774776114
$ws.Cells.Item(22, 7).Value = 3.775219841663274
$ws.Cells.Item(22, 9).Value = 59.25350203890755
$ws.Cells.Item(22, 10).Value = 8.371877872785397
$ws.Cells.Item(22, 13).Value = 25.17263998598118
$ws.Cells.Item(23, 2).Value = 22.16949047160742
$ws.Cells.Item(23, 3).Value = 26.77093083842269
$ws.Cells.Item(23, 4).Value = 14.00982581202729
$ws.Cells.Item(23, 5).Value = 14.07926352147798
$ws.Cells.Item(23, 7).Value = 3.777632110952569
$ws.Cells.Item(23, 9).Value = 58.69263987818394
$ws.Cells.Item(23, 10).Value = 8.375996479030952
$ws.Cells.Item(23, 13).Value = 25.14146247252059
$ws.Cells.Item(24, 2).Value = 21.94614418379164
$ws.Cells.Item(24, 3).Value = 25.9939789009049
$ws.Cells.Item(24, 4).Value = 13.96451677726742
$ws.Cells.Item(24, 5).Value = 14.07101196563313
$ws.Cells.Item(24, 7).Value = 3.787068902420473
$ws.Cells.Item(24, 9).Value = 56.52378916299797
$ws.Cells.Item(24, 10).Value = 8.392381024840653
$ws.Cells.Item(24, 13).Value = 25.03768110588201
$ws.Cells.Item(25, 2).Value = 21.74993803891655
$ws.Cells.Item(25, 3).Value = 25.16215271013741
$ws.Cells.Item(25, 4).Value = 13.92616216095064
$ws.Cells.Item(25, 5).Value = 14.06981829885579
$ws.Cells.Item(25, 7).Value = 3.797899960378158
$ws.Cells.Item(25, 9).Value = 54.09711053988136
$ws.Cells.Item(25, 10).Value = 8.41174524100874
$ws.Cells.Item(25, 13).Value = 24.95631999403772

Write-Host "updated 192 cells"